$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E) / Correspond Handback DateTime (H)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2:E3").Value = "2016-03-24 11:05:20"
$wsZh.Range("H2:H3").Value = "2016-03-24 11:06:01"

# de-de sheet: Correspond Handoff Datetime (E) / Correspond Handback DateTime (H)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2:E3").Value = "2016-03-24 11:05:25"
$wsDe.Range("H2:H3").Value = "2016-03-24 11:06:10"
